# Fixed excel reading for deposite testcase, like itself testcase
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet represented a "Deposite" test fixture, but was still named Sheet1.
$ws.Name = "Deposite"

# Header row: "sum" / "monthdep"
$ws.Range("A1").Value = "sum"
$ws.Range("B1").Value = "monthdep"

# Sample data row
$ws.Range("A2").Value = 1200
$ws.Range("B2").Value = 200

# Bump the font size across the new data range (also mints the style used
# by the header/data cells) and print the sheet in landscape.
$ws.Range("A1:B2").Font.Size = 11
$ws.PageSetup.Orientation = 2
